# Daily attendance processing - 2026-01-19 08:46:44
#
# For every row in the "Recorded By" column (G), when the comma-separated
# list of recorders contains the literal entry "System" (case-sensitive)
# but it is not already the first entry, move it to the front while
# preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    $parts = $val -split ', '

    if ($parts.Count -gt 1 -and -not ($parts[0].Equals("System"))) {
        $sysIndex = -1
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($parts[$i].Equals("System")) {
                $sysIndex = $i
            }
        }

        if ($sysIndex -ge 0) {
            $newParts = @("System")
            for ($i = 0; $i -lt $parts.Count; $i++) {
                if ($i -ne $sysIndex) {
                    $newParts += $parts[$i]
                }
            }
            $cell.Value2 = ($newParts -join ', ')
        }
    }
}
